$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a4 = @'
2025-12-28 01:13:43
'@
$b4 = @'
gemini-3-flash-preview
'@
$c4 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document f...
'@
$d4 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document fragment.
Input: The attached text from a contract.
CRITICAL INSTRUCTIONS:
1. **Assume Isolation with Common Sense**: Do NOT assume missing definitions exist in other documents. However, IGNORE common commercial lending terms typically defined in a base Credit Agreement (e.g., "Borrower", "Administrative Agent", "Lender", "Business Day", "Dollars", "GAAP", "Material Adverse Effect"). Only flag specific, deal-specific, or unusual capitalized terms that are undefined.
2. **Logic Check:** Check all math and logic tables.
3. **Drafting Errors:** Find any placeholders like "[__]" or blank lines that were forgotten.
Output Format:
Return ONLY a valid JSON object with the following structure:
{
  "errors": [
    {
      "location": "Page 3, Section 2.1",
      "error": "Description of the error",
      "suggestion": "Suggested fix"
    }
  ]
}
If no errors are found, return {"errors": []}.
--- CONTRACT TEXT BEGINS ---
--- [START OF PAGE 1] ---
FIRST AMENDMENT TO AMENDED AND RESTATED
CREDIT AGREEMENT
This document is a generated test file containing intentional legal drafting errors for AI training
purposes.
ARTICLE I: DEFINITIONS
...
"Applicable Margin" means the corresponding percentages per annum as set forth below based on
the Consolidated Total Leverage Ratio:
Pricing Level
Consolidated Total Leverage Ratio
SOFR Margin
I
Greater than or equal to 3.25 to 1.00
2.75%
II
Greater than or equal to 2.60 to 1.00 but less than 3.25 to 1.00
2.50%
III
Greater than or equal to 1.75 to 1.00 but less than 2.50 to 1.00
2.25%
IV
Greater than or equal to 1.50 to 1.00 but less than 1.75 to 1.00
2.00%
V
Less than 1.50 to 1.00
1.75%
"Cash Collateral" shall have a meaning correlative to the foregoing and shall include the proceeds of
such cash collateral...
(Note: "Cash Collateralize" is defined, but "Cash Collateralization" is NOT defined
in this section.)
SECTION 2.5 Permanent Reduction of the Revolving Credit
Commitment
The Borrower shall have the right at any time to terminate the Revolving Credit Commitment...
Any reduction of the Revolving Credit Commitment to zero shall be accompanied by payment of all
outstanding Revolving Credit Loans and furnishing of Cash Collateralization satisfactory to the
Administrative Agent.
SECTION 5.13 Incremental Loans
At any time after the First Amendment Effective Date, the Borrower may by written notice to the
Administrative Agent elect to request the establishment of one or more increases in the Revolving
--- [START OF PAGE 2] ---
Credit Commitments...
Such notice shall be delivered to the Administrative Agent at its office in [__] (or such other location as
the Administrative Agent may designate).
The terms of such Incremental Revolving Credit Increase shall be subject to the documentation
requirements set forth in Section 5.19.
[End of Document]
--- CONTRACT TEXT ENDS ---
'@
$e4 = @'
{
  "errors": [
    {
      "location": "Page 1, Article I, Definition of 'Applicable Margin'",
      "error": "Mathematical logic gap in pricing grid. The range between 2.50 to 1.00 and 2.60 to 1.00 is not covered. Level II starts at 2.60 while Level III ends at 2.50.",
      "suggestion": "Adjust the Consolidated Total Leverage Ratio for Pricing Level II to 'Greater than or equal to 2.50 to 1.00 but less than 3.25 to 1.00' to ensure continuity."
    },
    {
      "location": "Page 1, Section 2.5",
      "error": "The term 'Cash Collateralization' is used but is explicitly noted as being undefined in the document.",
      "suggestion": "Add a definition for 'Cash Collateralization' in Article I or replace the usage with the defined term 'Cash Collateral'."
    },
    {
      "location": "Page 2, Section 5.13",
      "error": "Drafting placeholder '[__]' found in the notice delivery provision.",
      "suggestion": "Replace '[__]' with the specific physical location or city of the Administrative Agent's office."
    },
    {
      "location": "Page 2, Section 5.13",
      "error": "Internal cross-reference to 'Section 5.19' cannot be verified as the section is missing from the document.",
      "suggestion": "Confirm that Section 5.19 exists in the full agreement or update the reference to the correct section."
    }
  ]
}
'@

$ws.Range("A4").Value = $a4
$ws.Range("B4").Value = $b4
$ws.Range("C4").Value = $c4
$ws.Range("D4").Value = $d4
$ws.Range("E4").Value = $e4

$ws.Rows.Item(4).EntireRow.AutoFit()

Write-Host "Row 4 written"
